# Backlog.xlsx update: rename existing sheet, add a new "TP2" sheet with
# its own backlog table, and select it as the active sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the existing sheet -----------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Product Backlog TP1"

# Clear the old selection/scroll position left over from the last TP1 edit.
$ws1.Range("A39:XFD40").Select()

# --- 2. Add the new "TP2" sheet right after it -----------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "TP2"

# --- 3. Header row -----------------------------------------------------
$ws2.Range("A1").Value = "Titre"
$ws2.Range("B1").Value = "Composant "
$ws2.Range("C1").Value = "Description"
$ws2.Range("D1").Value = "Membre"
$ws2.Range("E1").Value = "Status"

# --- 4. Backlog rows (row 2 left blank, data starts row 3) -----------------
$ws2.Range("A3").Value = "Watchdog "
$ws2.Range("B3").Value = "Watchdog"
$ws2.Range("C3").Value = "Implémentation du watchdog pour les 2 services core-backend"
$ws2.Range("D3").Value = "John"
$ws2.Range("E3").Value = "Fait"

$ws2.Range("A4").Value = "Cloud SQL "
$ws2.Range("B4").Value = "Base de données"
$ws2.Range("C4").Value = "Utilisation d'une base PostgreSQL directement depuis le cloud (GCP)"
$ws2.Range("D4").Value = "Williams"
$ws2.Range("E4").Value = "Fait"

$ws2.Range("A5").Value = "Photo associée au badge "
$ws2.Range("B5").Value = "Photo Storage Service"
$ws2.Range("C5").Value = "Création d'un service qui stocke les images et les met à disposition via une API (Request GET)"
$ws2.Range("D5").Value = "Williams"
$ws2.Range("E5").Value = "Fait"

$ws2.Range("A6").Value = "Replication Core-one/Core-Two "
$ws2.Range("B6").Value = "Core-one/Core-two"
$ws2.Range("C6").Value = "Réplication des services Core"
$ws2.Range("D6").Value = "John"
$ws2.Range("E6").Value = "Fait"

$ws2.Range("A7").Value = "Séparation Back-one"
$ws2.Range("B7").Value = "Back-one"
$ws2.Range("C7").Value = "Création d'un backend qui permet d'appeler Core-one/Core-two selon la disponibilité"
$ws2.Range("D7").Value = "John"
$ws2.Range("E7").Value = "Fait"

$ws2.Range("A8").Value = "Docker Swarm"
$ws2.Range("B8").Value = "CI/CD"
$ws2.Range("C8").Value = "Implémentation d'un système de CI/CD et de gestion des fichiers de config via Docker Swarm"
$ws2.Range("D8").Value = "Williams"
$ws2.Range("E8").Value = "Fait"

$ws2.Range("A9").Value = "Clusterisation Redis"
$ws2.Range("B9").Value = "Redis"
$ws2.Range("C9").Value = "Création d'un cluster Redis avec un Master et un Slave qui synchronisent les données"
$ws2.Range("D9").Value = "Ismail"
$ws2.Range("E9").Value = "Fait"

$ws2.Range("A10").Value = "Récupération des Images"
$ws2.Range("B10").Value = "Static-Server-Backend"
$ws2.Range("C10").Value = "Appel à l'API du service qui stocke les images et affichage sur le Front"
$ws2.Range("D10").Value = "Williams"
$ws2.Range("E10").Value = "Fait"

# --- 5. Column widths matching the source layout ---------------------------
$ws2.Columns.Item(1).ColumnWidth = 28.21875
$ws2.Columns.Item(2).ColumnWidth = 20.33203125
$ws2.Columns.Item(3).ColumnWidth = 77.44140625

# --- 6. Conditional formatting on the Status column (same palette as TP1) --
# À faire=yellow, En cours=gold, En revue=sandybrown, En test=skyblue,
# Bloqué=red, Fait=limegreen, Déployé=green (matches TP1's dxf colors).
$statusRange = $ws2.Range("E3:E10")
$rule1 = $statusRange.FormatConditions.Add(2, 3, '$E3="À faire"')
$rule1.Interior.Color = 65535
$rule2 = $statusRange.FormatConditions.Add(2, 3, '$E3="En cours"')
$rule2.Interior.Color = 55295
$rule3 = $statusRange.FormatConditions.Add(2, 3, '$E3="En revue"')
$rule3.Interior.Color = 6333684
$rule4 = $statusRange.FormatConditions.Add(2, 3, '$E3="En test"')
$rule4.Interior.Color = 15453831
$rule5 = $statusRange.FormatConditions.Add(2, 3, '$E3="Bloqué"')
$rule5.Interior.Color = 255
$rule6 = $statusRange.FormatConditions.Add(2, 3, '$E3="Fait"')
$rule6.Interior.Color = 3329330
$rule7 = $statusRange.FormatConditions.Add(2, 3, '$E3="Déployé"')
$rule7.Interior.Color = 32768

# --- 7. Selection / active sheet bookkeeping (must run last so the saved
#        file opens on TP2, matching activeTab/tabSelected in the source) --
$ws2.Range("E14").Select()
$ws2.Activate()
